# Re-order the "Pålmyrbäcken, Mpd" observation rows: a block of A62375-2025
# artfynd rows had their whole records (Id, activity, coordinates, public
# comment, ...) re-assigned to different row numbers - the underlying
# facts/rows are the same, only which worksheet row each record sits on
# changed.
#
# For each small cycle of row numbers (r0, r1, ..., rn-1) the new record
# shown on r[i] is the old record that used to be on r[i+1] (wrapping
# around) - i.e. every row "pulls" the record that used to sit one row
# later in the cycle.
#
# We touch every column that is actually populated on any of these rows,
# cell by cell, EXCEPT the Startdatum/Slutdatum columns (Y/AA): every row
# in every cycle carries the literal text "2026-01-25" there, so the
# values never actually change - and round-tripping a date-looking string
# through Range.Value gets reinterpreted as a real date serial, which
# would be an unwanted side effect for a column that isn't really being
# edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cycles of row numbers (1-based worksheet rows) whose records rotate.
$cycles = @(
    @(3, 4),
    @(9, 10, 11),
    @(21, 22, 23),
    @(24, 25),
    @(27, 28, 29),
    @(30, 33, 31, 32),
    @(35, 36),
    @(38, 39)
)

# Every column populated on the affected rows, excluding Y/AA (see above).
$cols = @("A","B","D","E","F","G","H","I","K","L","M","N","P","Q","R","S","T","U","V","W","AC","AD","AE","AG","AT","AW","AX","AY")

foreach ($cycle in $cycles) {
    $n = $cycle.Length

    # Snapshot every touched cell's old value before any writes happen.
    $old = @{}
    foreach ($r in $cycle) {
        foreach ($col in $cols) {
            $old["$col$r"] = $ws.Range("$col$r").Value()
        }
    }

    for ($i = 0; $i -lt $n; $i++) {
        $target = $cycle[$i]
        $source = $cycle[($i + 1) % $n]
        foreach ($col in $cols) {
            $ws.Range("$col$target").Value = $old["$col$source"]
        }
    }
}
